# Refresh cryptocurrency market data (price and 1h volume change %)
# Source: GitHub Actions scheduled data update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.535.25'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.60%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.861.70'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.85%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '520.08'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +6.53%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '141.13'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -3.81%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.608'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -2.18%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.15%  '

$ws.Range('E9').Value = '  -3.39%  '

$ws.Range('E10').Value = '  -6.11%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000319'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -7.41%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '41.56'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -3.47%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.30'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.475.73'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.75%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.46'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.54%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.858.00'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.20%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.14'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.71%  '

$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.133'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.81%  '

$ws.Range('B19').Value = 'Polygon'
$ws.Range('C19').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.20'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.79%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.569.86'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.52%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '416.61'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -4.69%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.47'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.98'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.20%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '86.73'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -3.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.99'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +6.91%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '11.41'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -8.60%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.51'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -5.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '35.39'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -4.84%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '13.18'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.64%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '677.04'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -4.82%  '

$ws.Range('E31').Value = '  -5.53%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.78'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -3.99%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.66'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +9.88%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '65.34'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +5.93%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.446'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -7.12%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '39.63'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.92%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0₃0829'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.84%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.51'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +13.91%  '

$ws.Range('E39').Value = '  -1.18%  '

$ws.Range('E40').Value = '  -0.06%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.12%  '

$ws.Range('E42').Value = '  -3.41%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.09'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.97%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.75'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -7.31%  '

$ws.Range('E45').Value = '  +1.48%  '

$ws.Range('E46').Value = '  -2.91%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.96'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -2.29%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000269'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +12.49%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '143.68'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.33%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.26'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.54%  '

$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0₆0334'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -9.21%  '
